$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab to reflect the new "through" date
$ws.Name = "Through 2022-05-19"

# Update the "May (through 05-18)" label to "May (through 05-19)"
$ws.Range("A6").Value = "May (through 05-19)"

# Update May row (row 6) values for years 2016-2019 and 2021-2022 (2020 unchanged)
$ws.Range("C6").Value = 31
$ws.Range("D6").Value = 38
$ws.Range("E6").Value = 26
$ws.Range("F6").Value = 27
$ws.Range("H6").Value = 71
$ws.Range("I6").Value = 65

# Update Total row (row 7) values for years 2016-2019 and 2021-2022 (2020 unchanged)
$ws.Range("C7").Value = 193
$ws.Range("D7").Value = 291
$ws.Range("E7").Value = 272
$ws.Range("F7").Value = 182
$ws.Range("H7").Value = 594
$ws.Range("I7").Value = 617
